{"js": "// The M2Doc AQL engine was upgraded from 7.x to 8.x; the wording of the\n// \"invalid expression\" diagnostic for an empty expression changed from\n// \"null or empty string.\" to \"missing expression\". Update the error text\n// that is rendered in the document accordingly (the surrounding\n// \"Invalid link statement: Expression \"\" is invalid: \" prefix is unchanged).\nconst oldText = \"null or empty string.\";\nconst newText = \"missing expression\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find the text to update: \"${oldText}\"`);\n}\n\n// Replace in place so the run is edited without disturbing any of the\n// surrounding runs/paragraphs (bookmarks, fields, other text runs).\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The M2Doc AQL engine was upgraded from 7.x to 8.x; the wording of the\n# \"invalid expression\" diagnostic for an empty expression changed from\n# \"null or empty string.\" to \"missing expression\". Update the error text\n# rendered in the document accordingly (the surrounding\n# 'Invalid link statement: Expression \"\" is invalid: ' prefix is unchanged).\n\n$d = $word.ActiveDocument\n\n$oldText = \"null or empty string.\"\n$newText = \"missing expression\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdReplaceAll = 2 : replace every match (there is exactly one in this document).\n$found = $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]2)\nif (-not $found) {\n    throw \"Could not find the text to update: '$oldText'\"\n}\n"}
